$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> (DAMSLTag, DialogAct) updates derived from the diff
$updates = @(
    @{ Row = 22; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 27; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 50; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 62; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 74; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 75; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 80; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 85; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 96; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 101; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 112; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 132; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 147; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 153; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 157; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 158; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 161; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 172; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 174; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 186; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 187; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 197; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 198; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 205; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 217; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 218; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 235; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 244; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 246; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 260; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 265; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 269; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 272; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 273; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 296; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 300; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 318; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 320; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 323; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 327; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 329; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 333; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 339; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 356; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 362; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 365; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 368; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 370; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 372; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 386; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 387; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 392; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 397; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 403; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 404; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 410; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 411; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 414; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 416; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 426; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 440; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 448; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 449; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 453; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 457; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 465; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 474; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

$wb.Save()